$wb = $excel.ActiveWorkbook

$wsVocab = $wb.Worksheets.Item("vocabulary")
$wsDesc  = $wb.Worksheets.Item("description")

# --- vocabulary sheet: insert new row 5 for "sample_date" attribute ---
$wsVocab.Rows.Item(5).Insert()
$wsVocab.Range("A5").Value = "sample_date"
$wsVocab.Range("C5").Value = "date"
$wsVocab.Range("D5").Value = "string"
$wsVocab.Range("E5").Value = "date of sampling format <dd/mm/yyyy>"

# --- description sheet: add new column I (level_05 header + "date" value for the category row) ---
$wsDesc.Range("I1").Value = "level_05"
$wsDesc.Range("I6").Value = "date"

# --- update legacy data validation on F4:F5 -> F4:F6, formula row 5 -> row 6 ---
$wsVocab.Range("F4:F6").Validation.Delete()
$wsVocab.Range("F4:F6").Validation.Add(3, 1, 1, "=`$G`$6:`$DB`$6")

# --- update (was in x14 extLst) validation for D2:D14 -> D2:D15 (unchanged formula) ---
$wsVocab.Range("D2:D15").Validation.Delete()
$wsVocab.Range("D2:D15").Validation.Add(3, 1, 1, "=description!`$E`$7:`$H`$7")

# --- replace C2:C14 validation (E6:H6) with a new C2:C1048576 validation (E6:I6) ---
$wsVocab.Range("C2:C14").Validation.Delete()
$wsVocab.Range("C2:C1048576").Validation.Add(3, 1, 1, "=description!`$E`$6:`$I`$6")

# --- stray formatting artifact: E21 gets date number format (numFmtId 14) + inherited wrap text ---
$wsVocab.Range("E21").NumberFormat = "mm-dd-yy"

# --- selection states (cosmetic, matches author's last click position) ---
$wsVocab.Range("H18").Select()
$wsDesc.Range("H9").Select()
$wsDesc.Application.ActiveWindow.ScrollColumn = 2
